# tests/header_select/has_header.xlsx — "add case for #167"
#
# The existing sample data (header row + 4 data rows) gets a new, blank
# row inserted in the middle (old row 3 "2, ok" is cleared out so there is
# a gap the header-detection code has to skip over), and a fresh data row
# is appended two rows below the last one (row 6 stays empty, row 7 holds
# the new case "6, ok"). This exercises the "has header" detector against
# a sheet whose used range contains an interior blank row as well as a
# trailing gap before the final row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-stamp the existing header + data block (A1:B5) with its own style.
# (Touching the explicit style here mirrors the workbook's own history —
# the block keeps the exact same font/number-format, just re-applied.)
$ws.Range("A1:B5").Style = "Normal"

# Row 3 (A3:B3) becomes a blank row in the middle of the data: drop its
# values but keep the cells/style in place.
$ws.Range("A3:B3").ClearContents()

# Leave row 6 empty and add the new test case two rows below the last
# existing row of data.
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "ok"

# Match the saved selection/active cell on the new last row.
$ws.Range("B7").Select() | Out-Null
